$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Wyposażenie" (Equipment) column (G) ---
$ws.Range("G1").Value = "Wyposażenie"
$ws.Range("G2").Value = "Projektor"
$ws.Range("G3").Value = "Projektor"
$ws.Range("G5").Value = "Projektor"
$ws.Range("G6").Value = "Projektor, Stanowiska komputerowe"
$ws.Range("G7").Value = "Projektor, Stanowiska komputerowe"
$ws.Range("G8").Value = "Stanowiska komputerowe, Urządzenia sieciowe i budowa sieci, Projektor"
$ws.Range("G9").Value = "Urządzenia sieciowe i budowa sieci, Projektor"
$ws.Range("G10").Value = "Projektor"
$ws.Range("G11").Value = "Projektor"
$ws.Range("G12").Value = "Projektor"
$ws.Range("G13").Value = "Oscyloskop, Płytka badawcza układu scalonego"
$ws.Range("G14").Value = "Projektor"
$ws.Range("G15").Value = "Płytka badawcza układu scalonego, Oscyloskop"
$ws.Range("G16").Value = "Płytka badawcza układu scalonego, Oscyloskop"

# --- Room-type rename for two computer labs (pracownia -> laboratorium) ---
$ws.Range("E8").Value = "laboratorium komputerowe"
$ws.Range("E12").Value = "laboratorium komputerowe"

# --- Column G width / formatting ---
$ws.Columns.Item(7).ColumnWidth = 66.5

# --- View state: scroll so column D is left-most visible, select G26 ---
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("G26").Select()
